# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values on Sheet1 for rows 2-40 with the newly
# regenerated strike-count values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New K values for rows 2..40 (row -> value), matching the regenerated
# save_data output.
$kValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 2
    6  = 2
    7  = 3
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 1
    13 = 2
    14 = 1
    15 = 0
    16 = 1
    17 = 1
    18 = 0
    19 = 1
    20 = 1
    21 = 3
    22 = 1
    23 = 6
    24 = 3
    25 = 2
    26 = 1
    27 = 2
    28 = 1
    29 = 2
    30 = 2
    31 = 5
    32 = 4
    33 = 0
    34 = 3
    35 = 2
    36 = 3
    37 = 2
    38 = 1
    39 = 1
    40 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
